# Generate Report for Handoff
# The "ed5744f3-..." localization entry moves from "Handed back: in sync
# with en-US" back to "Ready for handoff", and gets a fresh handoff
# timestamp recorded on the per-locale status sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = "Ready for handoff"
$zhcn.Range("D3").Value = "2016-03-08 16:47:17"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = "Ready for handoff"
$dede.Range("D3").Value = "2016-03-08 16:47:23"
